# "Final Changes in Front-End, updated Back-End"
#
# Updates the To-do/progress tracker on Tabelle1:
#  - Projektdokumentation section: Testplan and Installationsanleitung progress increased
#    (their summary % in C12 recalculates automatically from its formula)
#  - Back-End Applikation section: CRUD-Operationen, Validierung and Fehlerhandling are now
#    finished (100%) and flip from "in progress" (yellow) to "done" (green) status
#  - Back-End Vorbereitung section: Validierung finished (100%, yellow -> green) and
#    Unit-Tests progress increased
#  - the view scrolled down and the selection moved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status colors used by the legend/status column (B):
#   green  = Fertig (done)      -> RGB FF92D050 (BGR 5296274)
#   yellow = Angefangen (wip)   -> RGB FFFFFF00
#   red    = Nicht angefangen   -> RGB FFFF0000
$colorDone = 5296274

# --- Projektdokumentation ----------------------------------------------
# Testplan: 10 -> 30
$ws.Range("C9").Value = 30
# Installationsanleitung: 5 -> 10
$ws.Range("C10").Value = 10

# --- Back-End Applikation -----------------------------------------------
# CRUD-Operationen finished: 80 -> 100, status -> done
$ws.Range("C23").Value = 100
$ws.Range("B23").Interior.Color = $colorDone
# Validierung finished: 50 -> 100, status -> done
$ws.Range("C24").Value = 100
$ws.Range("B24").Interior.Color = $colorDone
# Fehlerhandling finished: 50 -> 100, status -> done
$ws.Range("C25").Value = 100
$ws.Range("B25").Interior.Color = $colorDone

# --- Back-End Vorbereitung -----------------------------------------------
# Validierung finished: 80 -> 100, status -> done
$ws.Range("C35").Value = 100
$ws.Range("B35").Interior.Color = $colorDone
# Unit-Tests: 5 -> 20
$ws.Range("C37").Value = 20

# --- View / selection ---------------------------------------------------
# Scroll the sheet down so row 22 is at the top and move the selection to
# cover D7 and E31 (E31 being the active cell).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D7,E31").Select() | Out-Null
$ws.Range("E31").Activate() | Out-Null
